{"js": "// Fill in the computed Retention and ARA (Answer Recall) values in the\n// single results table: each target row has a label in the first cell\n// and an (empty) value cell that receives a new bold run containing the\n// computed number.\nconst targets = {\n  \"Ratio\": \"0\",\n  \"Answer Recall Lenient (ARL)\": \"0.1666\",\n  \"Answer Recall Strict (ARS)\": \"0\",\n  \"Answer Recall Average (ARA)\": \"0.0833\",\n};\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const table = tables.items[t];\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < rows.items.length; i++) {\n    const row = rows.items[i];\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    if (cells.items.length !== 2) {\n      continue;\n    }\n\n    const labelCell = cells.items[0];\n    labelCell.load(\"value\");\n    await context.sync();\n\n    const label = (labelCell.value || \"\").trim();\n    if (!(label in targets)) {\n      continue;\n    }\n\n    const valueCell = cells.items[1];\n    const paragraphs = valueCell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    // The value cell has a single empty paragraph whose paragraph mark\n    // already carries bold / 12pt (sz 24 half-points) formatting; insert\n    // the computed value as a run there and stamp matching run formatting.\n    const para = paragraphs.items[paragraphs.items.length - 1];\n    const inserted = para.insertText(targets[label], \"Replace\");\n    inserted.font.bold = true;\n    inserted.font.size = 12;\n    inserted.font.sizeBidirectional = 12;\n    await context.sync();\n  }\n}\n", "ps1": "# Fill in the computed Retention and ARA (Answer Recall) values in the\n# single results table: each target row has a label in the first cell\n# and an (empty) value cell that receives a new bold run containing the\n# computed number.\n$targets = @{\n    \"Ratio\" = \"0\";\n    \"Answer Recall Lenient (ARL)\" = \"0.1666\";\n    \"Answer Recall Strict (ARS)\" = \"0\";\n    \"Answer Recall Average (ARA)\" = \"0.0833\";\n}\n\n$d = $word.ActiveDocument\n\nfor ($ti = 1; $ti -le $d.Tables.Count; $ti++) {\n    $table = $d.Tables.Item($ti)\n    for ($i = 1; $i -le $table.Rows.Count; $i++) {\n        $row = $table.Rows.Item($i)\n        if ($row.Cells.Count -ne 2) {\n            continue\n        }\n\n        $label = $row.Cells.Item(1).Range.Text\n        $label = $label.TrimEnd([char]7).TrimEnd([char]13).Trim()\n\n        if (-not $targets.ContainsKey($label)) {\n            continue\n        }\n\n        $valueRange = $row.Cells.Item(2).Range\n        # The value cell's single empty paragraph already carries bold /\n        # 12pt (sz 24 half-points) formatting on its paragraph mark; set\n        # the text then stamp matching run formatting so the inserted run\n        # carries it explicitly too.\n        $valueRange.Text = $targets[$label]\n        $valueRange.Font.Bold = $true\n        $valueRange.Font.Size = 12\n        $valueRange.Font.SizeBi = 12\n    }\n}\n"}
